$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("model_3_4_0", 0.9720713924244446, 0.9771709017512078, 0.9930971536222714, 0.9850971833511245, 3.368198388813411, 2.668521311884286, 0.7236873987402453, 1.753306086197026),
    @("model_3_4_1", 0.9758461888043121, 0.9734139189358506, 0.9924099350605413, 0.9828328310429932, 2.912956821528961, 3.107679644022705, 0.7957346943180477, 2.019705571386988),
    @("model_3_4_2", 0.978640357211946, 0.9691415738109713, 0.9912373053875952, 0.9800938200533248, 2.575979279658839, 3.607079309012464, 0.9186719974627237, 2.341948322639559),
    @("model_3_4_24", 0.9797180882405156, 0.8999481909089886, 0.963402239253654, 0.9320255443832532, 2.446004596739821, 11.69517875573851, 3.836872041610719, 7.997147757149975),
    @("model_3_4_23", 0.9799389045175277, 0.901197123696549, 0.9639002742969282, 0.9328912915064822, 2.419374087988401, 11.5491894694165, 3.784658553836826, 7.895293206025283),
    @("model_3_4_22", 0.9801774831464911, 0.9025717980659586, 0.9644505495748282, 0.9338451297102354, 2.390601434303252, 11.38850209527265, 3.72696825295755, 7.783074800116498),
    @("model_3_4_21", 0.9804341832031551, 0.9040832087615939, 0.9650580974466592, 0.93489488614779, 2.359643331042102, 11.21183144414318, 3.663273551537848, 7.659571680246338),
    @("model_3_4_3", 0.9806633333894496, 0.9645576514268105, 0.9897299369577317, 0.9770505903750362, 2.332007750349977, 4.14289962223065, 1.07670297168094, 2.699982192501436),
    @("model_3_4_20", 0.9807091010007269, 0.9057431226053771, 0.9657281066048691, 0.9360490183867199, 2.3264881628035, 11.01780207777249, 3.593030472334836, 7.523788819430917),
    @("model_3_4_19", 0.9810017697011846, 0.9075633435916092, 0.9664663756161949, 0.9373159654642241, 2.291192230391879, 10.80503421278, 3.515631099504994, 7.374733370772138),
    @("model_3_4_18", 0.9813112353516218, 0.9095563193721921, 0.9672785237355177, 0.9387048721553397, 2.253870580811747, 10.57207282785819, 3.430486316077566, 7.211329457802669),
    @("model_3_4_17", 0.9816356296926051, 0.911734520522727, 0.9681705596357868, 0.9402246995558815, 2.214748526707047, 10.31746022212021, 3.33696617888728, 7.032522813788843),
    @("model_3_4_16", 0.9819721287484311, 0.9141104887859921, 0.9691487282617324, 0.9418847223387977, 2.174166640388357, 10.03973037586055, 3.234415974275485, 6.837222279863775),
    @("model_3_4_4", 0.982083245852147, 0.9598193261305402, 0.9880082372249597, 0.9738362394719363, 2.16076588459725, 4.69676828134016, 1.257204221866889, 3.078148356278478),
    @("model_3_4_15", 0.9823162948286277, 0.9166959169108394, 0.9702184404599492, 0.9436931705868825, 2.132660108647872, 9.737516509315695, 3.122268434584329, 6.624459592470807),
    @("model_3_4_14", 0.9826621667312869, 0.919501958361235, 0.9713853342872206, 0.9456585077795236, 2.090947854210508, 9.409514880394259, 2.999932471670734, 6.393239029105211),
    @("model_3_4_13", 0.9830013616219287, 0.9225385626454219, 0.9726540558304115, 0.9477877469797841, 2.050040849410401, 9.054562478866952, 2.866920994509671, 6.142735507741176),
    @("model_3_4_5", 0.9830346873090676, 0.9550460234341603, 0.9861659246058715, 0.9705529408543965, 2.046021761619381, 5.254725491675365, 1.450350404472924, 3.464426171041054),
    @("model_3_4_12", 0.983322474627601, 0.925812585762406, 0.9740274360614071, 0.9500858385135009, 2.011314525321155, 8.671857898595752, 2.722937206885895, 5.872366625962877),
    @("model_3_4_11", 0.9836105290676123, 0.9293290826323126, 0.9755069346597023, 0.9525559296968136, 1.976575074093178, 8.260810263763451, 2.567828077484646, 5.581762104200168),
    @("model_3_4_6", 0.983625146344597, 0.950328004176317, 0.9842749618149134, 0.9672782200260093, 1.974812226136804, 5.806220553032042, 1.648597021653674, 3.849694814834025),
    @("model_3_4_10", 0.9838458229315418, 0.9330881587101711, 0.977091016134393, 0.9551975358349913, 1.948198564049089, 7.821407247602242, 2.401754585611758, 5.270978966453391),
    @("model_3_4_7", 0.9839403744263666, 0.9457319688983772, 0.9823887961653315, 0.9640697979649191, 1.936795625622344, 6.343456757269913, 1.846340711407972, 4.227163454443347),
    @("model_3_4_9", 0.9840027750903659, 0.9370847467115999, 0.9787744545343254, 0.9580056788476478, 1.929270086965429, 7.354241171202444, 2.225264614674956, 4.940602880441912),
    @("model_3_4_8", 0.9840483552532521, 0.941306541726081, 0.9805463641037134, 0.9609693796977168, 1.923773105750755, 6.860750370656017, 2.039499416247447, 4.591925522293946)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $rowData = $data[$i]
    $ws.Cells.Item($row, 1).Value = $rowData[0]
    for ($c = 1; $c -lt $rowData.Count; $c++) {
        $ws.Cells.Item($row, $c + 1).Value = $rowData[$c]
    }
}
